# Generate Report for Handoff
# - Mark the f3766040-1d41-419f-b0d8-3f001ac12be7.md file as "Ready for handoff"
#   (status was "Handed back: in sync with en-US") in the Overview sheet as well
#   as in both the zh-cn and de-de locale sheets.
# - Update the "Latest Handoff Datetime" for the b0887b31-035f-4d19-be4f-ea01448c9a65.md
#   file in both locale sheets (the f3766040 row shared the exact same handoff
#   timestamp string, so it is refreshed to the same new value too).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D2").Value = "2016-03-10 06:29:58"
$zhcn.Range("D3").Value = "2016-03-10 06:29:58"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D2").Value = "2016-03-10 06:30:09"
$dede.Range("D3").Value = "2016-03-10 06:30:09"
